# Update the "Förändrad" (Changed) date column (C) from 45523 to 45524
# for all data rows (2 through 28) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45523) {
        $cell.Value2 = 45524
    }
}
